# "Corrected the Shortest path finding algorithm."
#
# Sheet1 columns: A=Name, B=arena_length, C=arena_height, D=exit_size,
# E=slits, F=wallthick, G="used to be" (I/J/K mirror D/E/F for a second,
# "current" set of sizes).
#
# Semantic cell-value changes:
#   B2 (L  arena_length):  40   -> 30
#   D2 (L  exit_size):     3.7  -> 3.75
#   D3 (M  exit_size):     1.8  -> 1.85
#   F3 (M  wallthick):     0.2  -> 0.15
#   K3 (M  "used to be" wallthick, stored as TEXT): 0.25 -> 0.2
#   D4 (S  exit_size):     0.87 -> 0.92
#   B5 (XS arena_length):  4    -> 6
#   D5 (XS exit_size):     0.49 -> 0.54
#   E5 (XS slits, stored as TEXT): "2.6, 3.42" -> "2.6, 3.45"
#
# Plus a UI-state only change: active selection moves to L9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30
$ws.Range("D2").Value = 3.75

$ws.Range("D3").Value = 1.85
$ws.Range("F3").Value = 0.15

$ws.Range("D4").Value = 0.92

$ws.Range("B5").Value = 6
$ws.Range("D5").Value = 0.54
# "2.6, 3.45" is not numeric text, so a plain assignment keeps it a string.
$ws.Range("E5").Value = "2.6, 3.45"

# K3 must keep holding the text "0.2" (not the number 0.2). A direct
# Range.Value assignment of a numeric-looking string gets coerced back to
# a number by the engine, so instead: build "0.2" as a formula result
# (guaranteed text), paste only its *value* (as literal text) into a scratch
# cell, copy that into K3 (carries the text straight through), restore K3's
# original (fill-shaded) cell format from its neighbour J3, then wipe the
# scratch cells so nothing extra is left behind.
$ws.Range("Z1").Formula = "=""0.2"""
$ws.Range("Z1").Copy()
$ws.Range("Z2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z2").Copy($ws.Range("K3"))
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Z1:Z2").ClearContents()

$ws.Range("L9").Select()
